$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50; this shifts existing rows 50-57 down to 51-58
# and copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows.Item(50).Insert()

# Fill in the new row 50 with the new record's data.
$ws.Cells.Item(50, 1).Value = 6
$ws.Cells.Item(50, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(50, 3).Value = "Metropolitana"
$ws.Cells.Item(50, 4).Value = 44522
$ws.Cells.Item(50, 5).Value = 13
$ws.Cells.Item(50, 6).Value = 100114007
$ws.Cells.Item(50, 7).Value = "Jengibre"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 400
$ws.Cells.Item(50, 11).Value = 13000
$ws.Cells.Item(50, 12).Value = 15000
$ws.Cells.Item(50, 13).Value = 13850
$ws.Cells.Item(50, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(50, 15).Value = "Perú"
$ws.Cells.Item(50, 16).Value = 1065
$ws.Cells.Item(50, 17).Value = 13
$ws.Cells.Item(50, 18).Value = "Hortaliza"
